# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on row 5 of the
# zh-cn and de-de worksheets to reflect the regenerated report times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-27 08:44:40"
$wsZhCn.Range("G5").Value = "2016-01-27 08:45:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-27 08:44:52"
$wsDeDe.Range("G5").Value = "2016-01-27 08:45:47"
